$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: shift existing Acelga price rows down by two (a new
# week of data is inserted at the top of the date-ordered block) and
# append two more rows so the oldest week (44272) appears twice more
# at the bottom, matching the source export. Rows 294/295 receive the
# new week (44568).

# Row 294
$ws.Range("A294").Value = 8
$ws.Range("B294").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C294").Value = 'Coquimbo'
$ws.Range("D294").Value = 44568
$ws.Range("E294").Value = 4
$ws.Range("F294").Value = 100112009
$ws.Range("G294").Value = 'Acelga'
$ws.Range("H294").Value = 'Sin especificar'
$ws.Range("I294").Value = 'Primera'
$ws.Range("J294").Value = 2600
$ws.Range("K294").Value = 450
$ws.Range("L294").Value = 500
$ws.Range("M294").Value = 475
$ws.Range("N294").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O294").Value = 'Provincia del Elquí'
$ws.Range("P294").Value = 238
$ws.Range("Q294").Value = 2
$ws.Range("R294").Value = 'Hortaliza'

# Row 295
$ws.Range("A295").Value = 8
$ws.Range("B295").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C295").Value = 'Coquimbo'
$ws.Range("D295").Value = 44568
$ws.Range("E295").Value = 4
$ws.Range("F295").Value = 100112009
$ws.Range("G295").Value = 'Acelga'
$ws.Range("H295").Value = 'Sin especificar'
$ws.Range("I295").Value = 'Segunda'
$ws.Range("J295").Value = 1500
$ws.Range("K295").Value = 350
$ws.Range("L295").Value = 400
$ws.Range("M295").Value = 375
$ws.Range("N295").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O295").Value = 'Provincia del Elquí'
$ws.Range("P295").Value = 188
$ws.Range("Q295").Value = 2
$ws.Range("R295").Value = 'Hortaliza'

# Row 296
$ws.Range("A296").Value = 8
$ws.Range("B296").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C296").Value = 'Coquimbo'
$ws.Range("D296").Value = 44463
$ws.Range("E296").Value = 4
$ws.Range("F296").Value = 100112009
$ws.Range("G296").Value = 'Acelga'
$ws.Range("H296").Value = 'Sin especificar'
$ws.Range("I296").Value = 'Primera'
$ws.Range("J296").Value = 3400
$ws.Range("K296").Value = 450
$ws.Range("L296").Value = 500
$ws.Range("M296").Value = 475
$ws.Range("N296").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O296").Value = 'Provincia del Elquí'
$ws.Range("P296").Value = 238
$ws.Range("Q296").Value = 2
$ws.Range("R296").Value = 'Hortaliza'

# Row 297
$ws.Range("A297").Value = 8
$ws.Range("B297").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C297").Value = 'Coquimbo'
$ws.Range("D297").Value = 44463
$ws.Range("E297").Value = 4
$ws.Range("F297").Value = 100112009
$ws.Range("G297").Value = 'Acelga'
$ws.Range("H297").Value = 'Sin especificar'
$ws.Range("I297").Value = 'Segunda'
$ws.Range("J297").Value = 1500
$ws.Range("K297").Value = 350
$ws.Range("L297").Value = 400
$ws.Range("M297").Value = 375
$ws.Range("N297").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O297").Value = 'Provincia del Elquí'
$ws.Range("P297").Value = 188
$ws.Range("Q297").Value = 2
$ws.Range("R297").Value = 'Hortaliza'

# Row 298
$ws.Range("A298").Value = 8
$ws.Range("B298").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C298").Value = 'Coquimbo'
$ws.Range("D298").Value = 44365
$ws.Range("E298").Value = 4
$ws.Range("F298").Value = 100112009
$ws.Range("G298").Value = 'Acelga'
$ws.Range("H298").Value = 'Sin especificar'
$ws.Range("I298").Value = 'Primera'
$ws.Range("J298").Value = 3500
$ws.Range("K298").Value = 500
$ws.Range("L298").Value = 600
$ws.Range("M298").Value = 550
$ws.Range("N298").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O298").Value = 'Provincia del Elquí'
$ws.Range("P298").Value = 275
$ws.Range("Q298").Value = 2
$ws.Range("R298").Value = 'Hortaliza'

# Row 299
$ws.Range("A299").Value = 8
$ws.Range("B299").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C299").Value = 'Coquimbo'
$ws.Range("D299").Value = 44365
$ws.Range("E299").Value = 4
$ws.Range("F299").Value = 100112009
$ws.Range("G299").Value = 'Acelga'
$ws.Range("H299").Value = 'Sin especificar'
$ws.Range("I299").Value = 'Segunda'
$ws.Range("J299").Value = 1600
$ws.Range("K299").Value = 400
$ws.Range("L299").Value = 450
$ws.Range("M299").Value = 425
$ws.Range("N299").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O299").Value = 'Provincia del Elquí'
$ws.Range("P299").Value = 212
$ws.Range("Q299").Value = 2
$ws.Range("R299").Value = 'Hortaliza'

# Row 300
$ws.Range("A300").Value = 8
$ws.Range("B300").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C300").Value = 'Coquimbo'
$ws.Range("D300").Value = 44454
$ws.Range("E300").Value = 4
$ws.Range("F300").Value = 100112009
$ws.Range("G300").Value = 'Acelga'
$ws.Range("H300").Value = 'Sin especificar'
$ws.Range("I300").Value = 'Primera'
$ws.Range("J300").Value = 3400
$ws.Range("K300").Value = 450
$ws.Range("L300").Value = 500
$ws.Range("M300").Value = 475
$ws.Range("N300").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O300").Value = 'Provincia del Elquí'
$ws.Range("P300").Value = 238
$ws.Range("Q300").Value = 2
$ws.Range("R300").Value = 'Hortaliza'

# Row 301
$ws.Range("A301").Value = 8
$ws.Range("B301").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C301").Value = 'Coquimbo'
$ws.Range("D301").Value = 44454
$ws.Range("E301").Value = 4
$ws.Range("F301").Value = 100112009
$ws.Range("G301").Value = 'Acelga'
$ws.Range("H301").Value = 'Sin especificar'
$ws.Range("I301").Value = 'Segunda'
$ws.Range("J301").Value = 1600
$ws.Range("K301").Value = 350
$ws.Range("L301").Value = 400
$ws.Range("M301").Value = 375
$ws.Range("N301").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O301").Value = 'Provincia del Elquí'
$ws.Range("P301").Value = 188
$ws.Range("Q301").Value = 2
$ws.Range("R301").Value = 'Hortaliza'

# Row 302
$ws.Range("A302").Value = 8
$ws.Range("B302").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C302").Value = 'Coquimbo'
$ws.Range("D302").Value = 44272
$ws.Range("E302").Value = 4
$ws.Range("F302").Value = 100112009
$ws.Range("G302").Value = 'Acelga'
$ws.Range("H302").Value = 'Sin especificar'
$ws.Range("I302").Value = 'Primera'
$ws.Range("J302").Value = 3400
$ws.Range("K302").Value = 450
$ws.Range("L302").Value = 500
$ws.Range("M302").Value = 475
$ws.Range("N302").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O302").Value = 'Provincia del Elquí'
$ws.Range("P302").Value = 238
$ws.Range("Q302").Value = 2
$ws.Range("R302").Value = 'Hortaliza'

# Row 303
$ws.Range("A303").Value = 8
$ws.Range("B303").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C303").Value = 'Coquimbo'
$ws.Range("D303").Value = 44272
$ws.Range("E303").Value = 4
$ws.Range("F303").Value = 100112009
$ws.Range("G303").Value = 'Acelga'
$ws.Range("H303").Value = 'Sin especificar'
$ws.Range("I303").Value = 'Segunda'
$ws.Range("J303").Value = 1600
$ws.Range("K303").Value = 350
$ws.Range("L303").Value = 400
$ws.Range("M303").Value = 375
$ws.Range("N303").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O303").Value = 'Provincia del Elquí'
$ws.Range("P303").Value = 188
$ws.Range("Q303").Value = 2
$ws.Range("R303").Value = 'Hortaliza'

# Preserve/apply the date number format used elsewhere in column D
$ws.Range("D294:D303").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Dimension after edit:" $ws.UsedRange.Address()
